# Update the worksheet date and the 25 division-fact answers in the table.
# Cells are addressed directly by (row, col) via Table.Cell(...) so that
# duplicate/overlapping old-vs-new values (e.g. "42÷6=7, 0" and
# "83÷9=9, 2" appear both as a source value in one cell and as the target
# value of another) cannot cause a cascading find/replace mistake.

$d = $word.ActiveDocument

# Header date line.
$d.Content.Find.Execute("2026-01-31 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-01 Sunday", 2) | Out-Null

$t = $d.Tables.Item(1)

# Row 1 (table row 1)
$t.Cell(1, 1).Range.Text = "42÷6=7, 0"
$t.Cell(1, 2).Range.Text = "87÷8=10, 7"
$t.Cell(1, 3).Range.Text = "13÷7=1, 6"
$t.Cell(1, 4).Range.Text = "60÷9=6, 6"
$t.Cell(1, 5).Range.Text = "17÷8=2, 1"

# Row 2 (table row 5)
$t.Cell(5, 1).Range.Text = "31÷2=15, 1"
$t.Cell(5, 2).Range.Text = "29÷3=9, 2"
$t.Cell(5, 3).Range.Text = "99÷8=12, 3"
$t.Cell(5, 4).Range.Text = "54÷8=6, 6"
$t.Cell(5, 5).Range.Text = "31÷6=5, 1"

# Row 3 (table row 9)
$t.Cell(9, 1).Range.Text = "92÷4=23, 0"
$t.Cell(9, 2).Range.Text = "13÷9=1, 4"
$t.Cell(9, 3).Range.Text = "83÷9=9, 2"
$t.Cell(9, 4).Range.Text = "83÷3=27, 2"
$t.Cell(9, 5).Range.Text = "93÷6=15, 3"

# Row 4 (table row 13)
$t.Cell(13, 1).Range.Text = "90÷2=45, 0"
$t.Cell(13, 2).Range.Text = "12÷7=1, 5"
$t.Cell(13, 3).Range.Text = "16÷7=2, 2"
$t.Cell(13, 4).Range.Text = "24÷4=6, 0"
$t.Cell(13, 5).Range.Text = "61÷5=12, 1"

# Row 5 (table row 17)
$t.Cell(17, 1).Range.Text = "51÷6=8, 3"
$t.Cell(17, 2).Range.Text = "45÷3=15, 0"
$t.Cell(17, 3).Range.Text = "89÷7=12, 5"
$t.Cell(17, 4).Range.Text = "75÷5=15, 0"
$t.Cell(17, 5).Range.Text = "49÷5=9, 4"
